$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Лист1")

# Patient name fields (row 55)
$ws1.Range("C55").Value = "Желнов"
$ws1.Range("E55").Value = "Антон"
$ws1.Range("F55").Value = "Олегович"

# Date of birth fields (row 56): day, month, year -- these are text values, force text format
$ws1.Range("D56").NumberFormat = "@"
$ws1.Range("D56").Value = "3"
$ws1.Range("E56").NumberFormat = "@"
$ws1.Range("E56").Value = "2"

# Address fields (row 57): street, house number, apartment number -- force text format
$ws1.Range("C57").Value = "Лынькова"
$ws1.Range("E57").NumberFormat = "@"
$ws1.Range("E57").Value = "75"
$ws1.Range("G57").NumberFormat = "@"
$ws1.Range("G57").Value = "103"
